# Update the lattice-multiplication practice table: each of the 15 cells
# gets a new "AB x CD" problem (and the matching worked lines underneath)
# while the layout (title / top-factor line / dashes / two partial-product
# rows) stays the same. [char]11 is a manual line break (w:br) inside the
# cell's single run, matching the existing formatting.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-LatticeCell($table, $row, $col, $title, $topLine, $left, $right) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $title + [char]11 + $topLine + [char]11 + "  ----" + [char]11 + $left + [char]11 + $right
}

Set-LatticeCell $t 1 1 "26 x 47" "  4    7" "2|    |" "6|    |"
Set-LatticeCell $t 1 2 "12 x 52" "  5    2" "1|    |" "2|    |"
Set-LatticeCell $t 1 3 "90 x 46" "  4    6" "9|    |" "0|    |"

Set-LatticeCell $t 2 1 "21 x 61" "  6    1" "2|    |" "1|    |"
Set-LatticeCell $t 2 2 "88 x 66" "  6    6" "8|    |" "8|    |"
Set-LatticeCell $t 2 3 "40 x 76" "  7    6" "4|    |" "0|    |"

Set-LatticeCell $t 3 1 "43 x 58" "  5    8" "4|    |" "3|    |"
Set-LatticeCell $t 3 2 "51 x 23" "  2    3" "5|    |" "1|    |"
Set-LatticeCell $t 3 3 "59 x 44" "  4    4" "5|    |" "9|    |"

Set-LatticeCell $t 4 1 "58 x 88" "  8    8" "5|    |" "8|    |"
Set-LatticeCell $t 4 2 "34 x 85" "  8    5" "3|    |" "4|    |"
Set-LatticeCell $t 4 3 "11 x 86" "  8    6" "1|    |" "1|    |"

Set-LatticeCell $t 5 1 "58 x 37" "  3    7" "5|    |" "8|    |"
Set-LatticeCell $t 5 2 "67 x 69" "  6    9" "6|    |" "7|    |"
Set-LatticeCell $t 5 3 "80 x 80" "  8    0" "8|    |" "0|    |"

Write-Host "Lattice multiplication exercises updated."
